$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    0.0337525921132433,
    -0.00102444424694228,
    0.0736324550142593,
    -0.0493644672164725,
    -0.0411602433657146,
    -0.121544158018929,
    0.112519738309319,
    0.0283522521926804,
    0.0552313760906616,
    -0.0360267154150473,
    0.0182946147753705,
    0.0206824311319558,
    -0.0179937388482054,
    0.0276292582147641,
    0.0197943510483599,
    -0.0426506428566784,
    -0.000932168280232888,
    -0.0143321336842633,
    0.0635166560920207
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $col = $i + 1
    $ws.Cells.Item($row, 20).Value = $newValues[$i]
    $ws.Cells.Item(21, $col).Value = $newValues[$i]
}
